$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header shared-string runs (Volume number and week-covering dates)
$volChars = $ws.Range("A8").Characters(21, 2)
$volChars.Text = "25"

$weekRange = $ws.Range("C9")
$weekChars1 = $weekRange.Characters(27, 9)
$weekChars1.Text = "6/19/2023"
$weekChars2 = $weekRange.Characters(47, 9)
$weekChars2.Text = "6/25/2023"

# Update data table values (rows 14-30)
# Row 14
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = -20
$ws.Range("I14").Value = 33
$ws.Range("J14").Value = 33
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -23.255813953488
$ws.Range("M14").Value = -43.103448275862
$ws.Range("N14").Value = -84.722222222222

# Row 15
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = -40
$ws.Range("F15").Value = 14
$ws.Range("H15").Value = 7.692307692307
$ws.Range("I15").Value = 108
$ws.Range("J15").Value = 111
$ws.Range("K15").Value = -2.702702702702
$ws.Range("L15").Value = 3.846153846153
$ws.Range("M15").Value = 4.854368932038
$ws.Range("N15").Value = -63.636363636363

# Row 16
$ws.Range("C16").Value = 47
$ws.Range("D16").Value = 51
$ws.Range("E16").Value = -7.843137254901
$ws.Range("G16").Value = 202
$ws.Range("H16").Value = -7.425742574257
$ws.Range("I16").Value = 1136
$ws.Range("J16").Value = 1208
$ws.Range("K16").Value = -5.960264900662
$ws.Range("L16").Value = 23.478260869565
$ws.Range("M16").Value = -28.101265822784
$ws.Range("N16").Value = -85.277345775013

# Row 17
$ws.Range("C17").Value = 81
$ws.Range("D17").Value = 102
$ws.Range("E17").Value = -20.588235294117
$ws.Range("F17").Value = 334
$ws.Range("G17").Value = 385
$ws.Range("H17").Value = -13.246753246753
$ws.Range("I17").Value = 1995
$ws.Range("J17").Value = 1926
$ws.Range("K17").Value = 3.582554517133
$ws.Range("L17").Value = 27.232142857142
$ws.Range("M17").Value = 29.126213592233
$ws.Range("N17").Value = -49.468085106383

# Row 18
$ws.Range("C18").Value = 39
$ws.Range("D18").Value = 46
$ws.Range("E18").Value = -15.217391304347
$ws.Range("F18").Value = 136
$ws.Range("G18").Value = 184
$ws.Range("H18").Value = -26.086956521739
$ws.Range("I18").Value = 959
$ws.Range("J18").Value = 1173
$ws.Range("K18").Value = -18.243819266837
$ws.Range("L18").Value = 9.350057012542
$ws.Range("M18").Value = -29.329403095062
$ws.Range("N18").Value = -82.792033016328

# Row 19
$ws.Range("C19").Value = 118
$ws.Range("D19").Value = 128
$ws.Range("E19").Value = -7.8125
$ws.Range("F19").Value = 450
$ws.Range("G19").Value = 472
$ws.Range("H19").Value = -4.661016949152
$ws.Range("I19").Value = 2731
$ws.Range("J19").Value = 2691
$ws.Range("K19").Value = 1.486436269044
$ws.Range("L19").Value = 35.332011892963
$ws.Range("M19").Value = 49.072052401746
$ws.Range("N19").Value = -9.359442416196

# Row 20
$ws.Range("C20").Value = 30
$ws.Range("D20").Value = 26
$ws.Range("E20").Value = 15.384615384615
$ws.Range("F20").Value = 133
$ws.Range("G20").Value = 121
$ws.Range("H20").Value = 9.9173553719
$ws.Range("I20").Value = 821
$ws.Range("J20").Value = 813
$ws.Range("K20").Value = 0.984009840098
$ws.Range("L20").Value = 27.484472049689
$ws.Range("M20").Value = 21.091445427728
$ws.Range("N20").Value = -82.082060235704

# Row 21
$ws.Range("C21").Value = 320
$ws.Range("D21").Value = 359
$ws.Range("E21").Value = -10.863509749303
$ws.Range("F21").Value = 1258
$ws.Range("G21").Value = 1382
$ws.Range("H21").Value = -8.972503617945
$ws.Range("I21").Value = 7783
$ws.Range("J21").Value = 7955
$ws.Range("K21").Value = -2.162162162162
$ws.Range("L21").Value = 26.060900550696
$ws.Range("M21").Value = 8.807493359429
$ws.Range("N21").Value = -69.29177352535

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 25
$ws.Range("F22").Value = 13
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = -45.833333333333
$ws.Range("I22").Value = 138
$ws.Range("J22").Value = 176
$ws.Range("K22").Value = -21.590909090909
$ws.Range("L22").Value = 18.965517241379
$ws.Range("M22").Value = -31

# Row 23
$ws.Range("C23").Value = 33
$ws.Range("D23").Value = 29
$ws.Range("E23").Value = 13.793103448275
$ws.Range("F23").Value = 146
$ws.Range("G23").Value = 116
$ws.Range("H23").Value = 25.862068965517
$ws.Range("I23").Value = 780
$ws.Range("J23").Value = 706
$ws.Range("K23").Value = 10.481586402266
$ws.Range("L23").Value = 16.941529235382
$ws.Range("M23").Value = 48.854961832061

# Row 24
$ws.Range("C24").Value = 258
$ws.Range("D24").Value = 274
$ws.Range("E24").Value = -5.839416058394
$ws.Range("F24").Value = 983
$ws.Range("G24").Value = 1039
$ws.Range("H24").Value = -5.389797882579
$ws.Range("I24").Value = 5879
$ws.Range("J24").Value = 6167
$ws.Range("K24").Value = -4.670017836873
$ws.Range("L24").Value = 25.138356747552
$ws.Range("M24").Value = 25.727117194183

# Row 25
$ws.Range("C25").Value = 153
$ws.Range("D25").Value = 115
$ws.Range("E25").Value = 33.043478260869
$ws.Range("F25").Value = 540
$ws.Range("G25").Value = 575
$ws.Range("H25").Value = -6.086956521739
$ws.Range("I25").Value = 2922
$ws.Range("J25").Value = 2970
$ws.Range("K25").Value = -1.616161616161
$ws.Range("L25").Value = 41.775836972343
$ws.Range("M25").Value = -22.800528401585

# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -7.142857142857
$ws.Range("I26").Value = 164
$ws.Range("J26").Value = 179
$ws.Range("K26").Value = -8.379888268156
$ws.Range("L26").Value = -7.865168539325

# Row 27
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 11
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 60
$ws.Range("G27").Value = 54
$ws.Range("H27").Value = 11.111111111111
$ws.Range("I27").Value = 297
$ws.Range("J27").Value = 290
$ws.Range("K27").Value = 2.413793103448
$ws.Range("L27").Value = -8.615384615384

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 13
$ws.Range("E28").Value = -69.230769230769
$ws.Range("F28").Value = 17
$ws.Range("G28").Value = 40
$ws.Range("H28").Value = -57.5
$ws.Range("I28").Value = 110
$ws.Range("J28").Value = 147
$ws.Range("K28").Value = -25.17006802721
$ws.Range("L28").Value = -37.5
$ws.Range("M28").Value = -49.074074074074
$ws.Range("N28").Value = -87.885462555066

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = -50
$ws.Range("G29").Value = 29
$ws.Range("H29").Value = -44.827586206896
$ws.Range("I29").Value = 95
$ws.Range("J29").Value = 122
$ws.Range("K29").Value = -22.131147540983
$ws.Range("L29").Value = -37.5
$ws.Range("M29").Value = -44.444444444444
$ws.Range("N29").Value = -88.414634146341

# Row 30
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = -81.818181818181
$ws.Range("I30").Value = 31
$ws.Range("J30").Value = 39
$ws.Range("K30").Value = -20.51282051282
$ws.Range("L30").Value = -3.125
